$d = $word.ActiveDocument

function Merge-ParagraphRuns($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex).Range
    $p.Text = $newText
    $leftoverStart = $p.Start + $newText.Length
    $leftoverEnd = $p.End - 1
    if ($leftoverEnd -gt $leftoverStart) {
        $leftover = $d.Range($leftoverStart, $leftoverEnd)
        $leftover.Delete()
    }
}

Merge-ParagraphRuns 1 "Answers: Arithmetic on complex numbers"
Merge-ParagraphRuns 2 "Charlotte McCarthy"
Merge-ParagraphRuns 4 "Answers to questions relating to the guide on arithmetic on complex numbers."
